$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 31 values (week of 6/6/2025)
$ws.Range("D31").Value = (Get-Date -Year 2025 -Month 6 -Day 6).Date
$ws.Range("E31").Value = 100
$ws.Range("F31").Value = 261
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 1012
$ws.Range("J31").Value = "Seguire trabajando en fin de semana (rafael)"

# Update the view/selection state to match what was saved
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.TopLeftCell = $ws.Range("A15")
$ws.Range("D32").Select()
